$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 494.25
$ws.Range("I41").Value = 389
$ws.Range("J41").Value = 599.5
$ws.Range("K41").Value = 389
$ws.Range("L41").Value = 599.5
$ws.Range("M41").Value = 51
$ws.Range("N41").Value = -1479.5
$ws.Range("H51").Value = 12713.286
$ws.Range("I51").Value = 9664.333000000001
$ws.Range("K51").Value = 9664.333000000001
$ws.Range("M51").Value = -9180.333000000001
$ws.Range("H74").Value = 17583.166
$ws.Range("I74").Value = 18272.545
$ws.Range("K74").Value = 18272.545
$ws.Range("M74").Value = -17336.545
$ws.Range("H76").Value = 3250
$ws.Range("I76").Value = 3500
$ws.Range("K76").Value = 3500
$ws.Range("M76").Value = -3185
$ws.Range("H77").Value = 17583.166
$ws.Range("I77").Value = 18272.545
$ws.Range("K77").Value = 91362.72499999999
$ws.Range("M77").Value = -86682.72499999999
$ws.Range("H79").Value = 3250
$ws.Range("I79").Value = 3500
$ws.Range("K79").Value = 3500
$ws.Range("M79").Value = -2408
$ws.Range("H86").Value = 8463
$ws.Range("I86").Value = 8509.714
$ws.Range("K86").Value = 8509.714
$ws.Range("M86").Value = -7386.714
$ws.Range("H89").Value = 8463
$ws.Range("I89").Value = 8509.714
$ws.Range("K89").Value = 42548.57
$ws.Range("M89").Value = -36932.57
$ws.Range("H97").Value = 2250
$ws.Range("J97").Value = 2250
$ws.Range("L97").Value = 6750
$ws.Range("N97").Value = -7742
$ws.Range("H125").Value = 987.63635
$ws.Range("I125").Value = 954.6667
$ws.Range("K125").Value = 8592.0003
$ws.Range("M125").Value = -6132.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3212.8572
$ws.Range("I61").Value = 2622.75
$ws.Range("K61").Value = 2622.75
$ws.Range("M61").Value = -2410.75
$ws.Range("H88").Value = 2675.7856
$ws.Range("J88").Value = 3297.125
$ws.Range("L88").Value = 3297.125
$ws.Range("N88").Value = -4109.125
$ws.Range("H91").Value = 2675.7856
$ws.Range("J91").Value = 3297.125
$ws.Range("L91").Value = 3297.125
$ws.Range("N91").Value = -6105.125
$ws.Range("H101").Value = 32250
$ws.Range("J101").Value = 32250
$ws.Range("L101").Value = 32250
$ws.Range("N101").Value = -38740
$ws.Range("H132").Value = 2104.3794
$ws.Range("I132").Value = 1712.6666
$ws.Range("J132").Value = 2745.3635
$ws.Range("K132").Value = 5137.9998
$ws.Range("L132").Value = 8236.0905
$ws.Range("M132").Value = -2607.9998
$ws.Range("N132").Value = -13296.0905
$ws.Range("H136").Value = 3212.8572
$ws.Range("I136").Value = 2622.75
$ws.Range("K136").Value = 7868.25
$ws.Range("M136").Value = -5318.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 333350530
$ws.Range("I94").Value = 400020220
$ws.Range("K94").Value = 400020220
$ws.Range("M94").Value = -400019769
$ws.Range("I105").Value = 834904.7
$ws.Range("K105").Value = 834904.7
$ws.Range("M105").Value = -833157.7
$ws.Range("H134").Value = 2375.9167
$ws.Range("I134").Value = 1626.375
$ws.Range("K134").Value = 4879.125
$ws.Range("M134").Value = -2344.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7355285.5
$ws.Range("I31").Value = 2305.9285
$ws.Range("J31").Value = 41669190
$ws.Range("K31").Value = 2305.9285
$ws.Range("L31").Value = 41669190
$ws.Range("M31").Value = -2010.9285
$ws.Range("N31").Value = -41669780
$ws.Range("H34").Value = 7355285.5
$ws.Range("I34").Value = 2305.9285
$ws.Range("J34").Value = 41669190
$ws.Range("K34").Value = 2305.9285
$ws.Range("L34").Value = 41669190
$ws.Range("M34").Value = -2103.9285
$ws.Range("N34").Value = -41669594
$ws.Range("H58").Value = 1706.2106
$ws.Range("I58").Value = 1172.9286
$ws.Range("J58").Value = 3199.4
$ws.Range("K58").Value = 1172.9286
$ws.Range("L58").Value = 3199.4
$ws.Range("M58").Value = -969.9286
$ws.Range("N58").Value = -3605.4
$ws.Range("H86").Value = 5361.9375
$ws.Range("J86").Value = 6008
$ws.Range("L86").Value = 6008
$ws.Range("N86").Value = -8254
$ws.Range("H89").Value = 5361.9375
$ws.Range("J89").Value = 6008
$ws.Range("L89").Value = 30040
$ws.Range("N89").Value = -41272
$ws.Range("H132").Value = 4056.697
$ws.Range("I132").Value = 3505.8096
$ws.Range("K132").Value = 10517.4288
$ws.Range("M132").Value = -7987.4288
$ws.Range("H136").Value = 1706.2106
$ws.Range("I136").Value = 1172.9286
$ws.Range("J136").Value = 3199.4
$ws.Range("K136").Value = 3518.7858
$ws.Range("L136").Value = 9598.200000000001
$ws.Range("M136").Value = -968.7857999999997
$ws.Range("N136").Value = -14698.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39.5
$ws.Range("J12").Value = 49.125
$ws.Range("L12").Value = 147.375
$ws.Range("N12").Value = -493.375
$ws.Range("H33").Value = 271
$ws.Range("I33").Value = 270.63635
$ws.Range("J33").Value = 271.8
$ws.Range("K33").Value = 1623.8181
$ws.Range("L33").Value = 1630.8
$ws.Range("M33").Value = -1340.8181
$ws.Range("N33").Value = -2196.8
$ws.Range("H113").Value = 1203.2
$ws.Range("J113").Value = 1327.5714
$ws.Range("L113").Value = 3982.7142
$ws.Range("N113").Value = -8322.7142
$ws.Range("H117").Value = 1271.1666
$ws.Range("J117").Value = 1249.75
$ws.Range("L117").Value = 3749.25
$ws.Range("N117").Value = -10633.25
$ws.Range("H121").Value = 5617090
$ws.Range("I121").Value = 12500427
$ws.Range("K121").Value = 37501281
$ws.Range("M121").Value = -37499971
$ws.Range("H122").Value = 1059.6666
$ws.Range("I122").Value = 585.6667
$ws.Range("J122").Value = 1533.6666
$ws.Range("K122").Value = 5271.0003
$ws.Range("L122").Value = 13802.9994
$ws.Range("M122").Value = -2821.0003
$ws.Range("N122").Value = -18702.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 66669090
$ws.Range("I80").Value = 125002210
$ws.Range("K80").Value = 125002210
$ws.Range("M80").Value = -125001212
$ws.Range("H83").Value = 66669090
$ws.Range("I83").Value = 125002210
$ws.Range("K83").Value = 625011050
$ws.Range("M83").Value = -625006058
$ws.Range("H107").Value = 1549.1875
$ws.Range("J107").Value = 2446.2222
$ws.Range("L107").Value = 2446.2222
$ws.Range("N107").Value = -6286.2222
$ws.Range("H122").Value = 3638.1738
$ws.Range("I122").Value = 3110.8
$ws.Range("K122").Value = 9332.400000000001
$ws.Range("M122").Value = -6882.400000000001
$ws.Range("H132").Value = 2385.7407
$ws.Range("I132").Value = 2339.7222
$ws.Range("J132").Value = 2477.7778
$ws.Range("K132").Value = 7019.1666
$ws.Range("L132").Value = 7433.3334
$ws.Range("M132").Value = -4489.1666
$ws.Range("N132").Value = -12493.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2992.4546
$ws.Range("I7").Value = 2866.125
$ws.Range("J7").Value = 3329.3333
$ws.Range("K7").Value = 2866.125
$ws.Range("L7").Value = 3329.3333
$ws.Range("M7").Value = -2754.125
$ws.Range("N7").Value = -3553.3333
$ws.Range("H16").Value = 1541.7142
$ws.Range("I16").Value = 1465.3334
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1465.3334
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1295.3334
$ws.Range("N16").Value = -2340
$ws.Range("H126").Value = 2992.4546
$ws.Range("I126").Value = 2866.125
$ws.Range("J126").Value = 3329.3333
$ws.Range("K126").Value = 8598.375
$ws.Range("L126").Value = 9987.999899999999
$ws.Range("M126").Value = -6128.375
$ws.Range("N126").Value = -14927.9999
$ws.Range("H136").Value = 5051.6313
$ws.Range("I136").Value = 5004.5293
$ws.Range("J136").Value = 5452
$ws.Range("K136").Value = 15013.5879
$ws.Range("L136").Value = 16356
$ws.Range("M136").Value = -12463.5879
$ws.Range("N136").Value = -21456
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 43275
$ws.Range("J54").Value = 43275
$ws.Range("L54").Value = 43275
$ws.Range("N54").Value = -44315
$ws.Range("H126").Value = 11667.25
$ws.Range("I126").Value = 14223.111
$ws.Range("J126").Value = 3999.6667
$ws.Range("K126").Value = 42669.333
$ws.Range("L126").Value = 11999.0001
$ws.Range("M126").Value = -40199.333
$ws.Range("N126").Value = -16939.0001
